$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.051.78"
$ws.Range("E2").Value = "  -2.06%  "
$ws.Range("D3").Value = "2.911.43"
$ws.Range("E3").Value = "  -2.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "371.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.01%  "
$ws.Range("E7").Value = "  -3.64%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.584"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.02%  "
$ws.Range("E10").Value = "  -3.80%  "
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("E12").Value = "  -2.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.86%  "
$ws.Range("D14").Value = "3.369.90"
$ws.Range("E14").Value = "  -2.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.58%  "
$ws.Range("D16").Value = "2.913.68"
$ws.Range("E16").Value = "  -1.98%  "
$ws.Range("E17").Value = "  -8.88%  "
$ws.Range("D18").Value = "50.989.30"
$ws.Range("E18").Value = "  -2.26%  "
$ws.Range("E19").Value = "  -7.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.23%  "
$ws.Range("D22").Value = "0.0₃0940"
$ws.Range("E22").Value = "  -3.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "258.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.36%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -6.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.101"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.40%  "
$ws.Range("E33").Value = "  -2.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("E35").Value = "  -5.79%  "
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0420"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.46%  "
$ws.Range("E38").Value = "  -7.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.09%  "
$ws.Range("E41").Value = "  -6.43%  "
$ws.Range("E42").Value = "  -3.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "119.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.68%  "
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("D46").Value = "2.019.22"
$ws.Range("E46").Value = "  -4.67%  "
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("E48").Value = "  -7.34%  "
$ws.Range("D49").Value = "3.197.94"
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.234"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("E51").Value = "  -8.90%  "
